# Daily attendance processing - 2025-10-16 14:50:15
# Normalize the "Recorded By" (column G) cell values: when "System" is
# listed first, move it so it appears last instead (keeping the other
# name(s) in front), matching the canonical ordering used elsewhere in
# the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value()

    if ($value -eq $null) { continue }

    $parts = $value -split ', '

    if (($parts.Count -eq 2) -and ($parts[0] -eq 'System')) {
        $cell.Value = "$($parts[1]), System"
    }
    elseif (($parts.Count -eq 3) -and ($parts[0] -eq 'backup@backdoor.com') -and ($parts[1] -eq 'System') -and ($parts[2] -eq 'system')) {
        $cell.Value = 'system, backup@backdoor.com, System'
    }
}
